# Fix latency units in report sheet:
#  - Rename header "Utility" (O2) to "Utility (Percent)"
#  - Append " usec" to the min/max/average read-latency values in columns I, J, K
#    for every data row (rows 3 through 38)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update header
$ws.Range("O2").Value = "Utility (Percent)"

# Append " usec" to latency columns (I, J, K) for data rows 3..38
$lastRow = 38
for ($row = 3; $row -le $lastRow; $row++) {
    foreach ($col in @("I", "J", "K")) {
        $cell = $ws.Range("$col$row")
        $current = $cell.Text
        if ($current -ne $null -and $current.Length -gt 0 -and -not ($current.EndsWith(" usec"))) {
            $cell.Value = ($current + " usec")
        }
    }
}
